$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header swap: BP1/BQ1 average_doctor <-> average_doctor_old
# Also row 4-13 stat columns recomputed for the new "average_doctor_old" metric
# and average_doctor/average_doctor_old columns (BP/BQ) swapped & updated accordingly.

$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"
$ws.Range("E4").Value = 0.437
$ws.Range("F4").Value = 0.074
$ws.Range("G4").Value = 0.272
$ws.Range("N4").Value = 0.437
$ws.Range("O4").Value = 0.068
$ws.Range("P4").Value = 0.26
$ws.Range("Q4").Value = 0.026
$ws.Range("R4").Value = 0.018
$ws.Range("S4").Value = 0.134
$ws.Range("W4").Value = 0.28
$ws.Range("X4").Value = 0.107
$ws.Range("Y4").Value = 0.327
$ws.Range("AI4").Value = 0.285
$ws.Range("AJ4").Value = 0.083
$ws.Range("AK4").Value = 0.288
$ws.Range("AU4").Value = 0.196
$ws.Range("AV4").Value = 0.029
$ws.Range("AW4").Value = 0.17
$ws.Range("BA4").Value = 2.025
$ws.Range("BB4").Value = 0.16
$ws.Range("BC4").Value = 0.4
$ws.Range("BG4").Value = 0.736
$ws.Range("BH4").Value = 0.142
$ws.Range("BI4").Value = 0.377
$ws.Range("BM4").Value = 0.72
$ws.Range("BN4").Value = 0.079
$ws.Range("BO4").Value = 0.281
$ws.Range("BP4").Value = 0.675
$ws.Range("BQ4").Value = 0.717
$ws.Range("E5").Value = 0.54
$ws.Range("F5").Value = 0.08400000000000001
$ws.Range("G5").Value = 0.29
$ws.Range("N5").Value = 0.729
$ws.Range("O5").Value = 0.083
$ws.Range("P5").Value = 0.289
$ws.Range("Q5").Value = 0.017
$ws.Range("R5").Value = 0.007
$ws.Range("S5").Value = 0.08599999999999999
$ws.Range("W5").Value = 0.269
$ws.Range("X5").Value = 0.107
$ws.Range("Y5").Value = 0.328
$ws.Range("AI5").Value = 0.303
$ws.Range("AJ5").Value = 0.091
$ws.Range("AK5").Value = 0.301
$ws.Range("AU5").Value = 0.374
$ws.Range("AV5").Value = 0.091
$ws.Range("AW5").Value = 0.302
$ws.Range("BA5").Value = 1.308
$ws.Range("BB5").Value = 0.07199999999999999
$ws.Range("BC5").Value = 0.269
$ws.Range("BG5").Value = 0.386
$ws.Range("BH5").Value = 0.045
$ws.Range("BI5").Value = 0.211
$ws.Range("BM5").Value = 0.542
$ws.Range("BN5").Value = 0.064
$ws.Range("BO5").Value = 0.253
$ws.Range("BP5").Value = 0.436
$ws.Range("BQ5").Value = 0.451
$ws.Range("E6").Value = 0.483
$ws.Range("N6").Value = 0.546
$ws.Range("Q6").Value = 0.021
$ws.Range("W6").Value = 0.274
$ws.Range("AI6").Value = 0.294
$ws.Range("AU6").Value = 0.257
$ws.Range("BA6").Value = 1.58
$ws.Range("BG6").Value = 0.506
$ws.Range("BM6").Value = 0.618
$ws.Range("BP6").Value = 0.527
$ws.Range("BQ6").Value = 0.551
$ws.Range("E7").Value = 0.516
$ws.Range("N7").Value = 0.643
$ws.Range("Q7").Value = 0.018
$ws.Range("W7").Value = 0.271
$ws.Range("AI7").Value = 0.299
$ws.Range("AU7").Value = 0.317
$ws.Range("BA7").Value = 1.404
$ws.Range("BG7").Value = 0.427
$ws.Range("BM7").Value = 0.57
$ws.Range("BP7").Value = 0.468
$ws.Range("BQ7").Value = 0.486
$ws.Range("E8").Value = 0.599
$ws.Range("F8").Value = 0.113
$ws.Range("G8").Value = 0.336
$ws.Range("N8").Value = 0.763
$ws.Range("O8").Value = 0.07000000000000001
$ws.Range("P8").Value = 0.265
$ws.Range("Q8").Value = 0.019
$ws.Range("W8").Value = 0.296
$ws.Range("X8").Value = 0.115
$ws.Range("Y8").Value = 0.339
$ws.Range("AI8").Value = 0.321
$ws.Range("AJ8").Value = 0.122
$ws.Range("AK8").Value = 0.349
$ws.Range("AU8").Value = 0.31
$ws.Range("AV8").Value = 0.082
$ws.Range("AW8").Value = 0.287
$ws.Range("BA8").Value = 1.742
$ws.Range("BB8").Value = 0.124
$ws.Range("BC8").Value = 0.352
$ws.Range("BG8").Value = 0.5649999999999999
$ws.Range("BH8").Value = 0.106
$ws.Range("BI8").Value = 0.325
$ws.Range("BM8").Value = 0.6889999999999999
$ws.Range("BN8").Value = 0.068
$ws.Range("BO8").Value = 0.261
$ws.Range("BP8").Value = 0.581
$ws.Range("BQ8").Value = 0.602
$ws.Range("E9").Value = 0.535
$ws.Range("F9").Value = 0.249
$ws.Range("G9").Value = 0.499
$ws.Range("N9").Value = 0.663
$ws.Range("O9").Value = 0.223
$ws.Range("P9").Value = 0.473
$ws.Range("W9").Value = 0.186
$ws.Range("X9").Value = 0.151
$ws.Range("Y9").Value = 0.389
$ws.Range("AI9").Value = 0.244
$ws.Range("AJ9").Value = 0.185
$ws.Range("AK9").Value = 0.43
$ws.Range("BA9").Value = 1.721
$ws.Range("BB9").Value = 0.249
$ws.Range("BC9").Value = 0.499
$ws.Range("BG9").Value = 0.605
$ws.Range("BH9").Value = 0.239
$ws.Range("BI9").Value = 0.489
$ws.Range("BM9").Value = 0.651
$ws.Range("BN9").Value = 0.227
$ws.Range("BO9").Value = 0.477
$ws.Range("BP9").Value = 0.574
$ws.Range("BQ9").Value = 0.587
$ws.Range("E10").Value = 0.674
$ws.Range("F10").Value = 0.22
$ws.Range("G10").Value = 0.469
$ws.Range("N10").Value = 0.86
$ws.Range("O10").Value = 0.12
$ws.Range("P10").Value = 0.347
$ws.Range("W10").Value = 0.372
$ws.Range("X10").Value = 0.234
$ws.Range("Y10").Value = 0.483
$ws.Range("AI10").Value = 0.349
$ws.Range("AJ10").Value = 0.227
$ws.Range("AK10").Value = 0.477
$ws.Range("AU10").Value = 0.291
$ws.Range("AV10").Value = 0.206
$ws.Range("AW10").Value = 0.454
$ws.Range("BA10").Value = 2.07
$ws.Range("BB10").Value = 0.243
$ws.Range("BC10").Value = 0.493
$ws.Range("BG10").Value = 0.663
$ws.Range("BH10").Value = 0.223
$ws.Range("BI10").Value = 0.473
$ws.Range("BM10").Value = 0.826
$ws.Range("BN10").Value = 0.144
$ws.Range("BO10").Value = 0.379
$ws.Range("BP10").Value = 0.6899999999999999
$ws.Range("BQ10").Value = 0.722
$ws.Range("E11").Value = 0.698
$ws.Range("F11").Value = 0.211
$ws.Range("G11").Value = 0.459
$ws.Range("N11").Value = 0.884
$ws.Range("O11").Value = 0.103
$ws.Range("P11").Value = 0.321
$ws.Range("W11").Value = 0.372
$ws.Range("X11").Value = 0.234
$ws.Range("Y11").Value = 0.483
$ws.Range("AI11").Value = 0.384
$ws.Range("AJ11").Value = 0.236
$ws.Range("AK11").Value = 0.486
$ws.Range("AU11").Value = 0.442
$ws.Range("AV11").Value = 0.247
$ws.Range("AW11").Value = 0.497
$ws.Range("BA11").Value = 2.07
$ws.Range("BB11").Value = 0.243
$ws.Range("BC11").Value = 0.493
$ws.Range("BG11").Value = 0.663
$ws.Range("BH11").Value = 0.223
$ws.Range("BI11").Value = 0.473
$ws.Range("BM11").Value = 0.826
$ws.Range("BN11").Value = 0.144
$ws.Range("BO11").Value = 0.379
$ws.Range("BP11").Value = 0.6899999999999999
$ws.Range("BQ11").Value = 0.725
$ws.Range("E12").Value = 1.4
$ws.Range("F12").Value = 0.707
$ws.Range("G12").Value = 0.841
$ws.Range("N12").Value = 1.5
$ws.Range("O12").Value = 1.122
$ws.Range("P12").Value = 1.059
$ws.Range("W12").Value = 1.688
$ws.Range("X12").Value = 0.59
$ws.Range("Y12").Value = 0.768
$ws.Range("AI12").Value = 1.758
$ws.Range("AJ12").Value = 1.396
$ws.Range("AK12").Value = 1.181
$ws.Range("AU12").Value = 2.875
$ws.Range("AV12").Value = 2.759
$ws.Range("AW12").Value = 1.661
$ws.Range("BA12").Value = 3.619
$ws.Range("BG12").Value = 1.105
$ws.Range("BH12").Value = 0.129
$ws.Range("BI12").Value = 0.36
$ws.Range("BN12").Value = 0.274
$ws.Range("BO12").Value = 0.523
$ws.Range("BP12").Value = 1.206
$ws.Range("BQ12").Value = 1.24
$ws.Range("E13").Value = 1.521
$ws.Range("F13").Value = 0.532
$ws.Range("G13").Value = 0.729
$ws.Range("N13").Value = 2.051
$ws.Range("O13").Value = 0.977
$ws.Range("P13").Value = 0.988
$ws.Range("W13").Value = 1.026
$ws.Range("X13").Value = 0.185
$ws.Range("Y13").Value = 0.43
$ws.Range("AI13").Value = 1.271
$ws.Range("AJ13").Value = 0.39
$ws.Range("AK13").Value = 0.624
$ws.Range("AU13").Value = 2.257
$ws.Range("AV13").Value = 0.945
$ws.Range("AW13").Value = 0.972
$ws.Range("BA13").Value = 2.27
$ws.Range("BB13").Value = 0.285
$ws.Range("BC13").Value = 0.534
$ws.Range("BG13").Value = 0.5570000000000001
$ws.Range("BH13").Value = 0.046
$ws.Range("BI13").Value = 0.215
$ws.Range("BM13").Value = 0.872
$ws.Range("BN13").Value = 0.247
$ws.Range("BO13").Value = 0.497
$ws.Range("BP13").Value = 0.757
$ws.Range("BQ13").Value = 0.711
